# Update Zalera_Profits market-data values (currentAveragePrice* / Leve* columns)
# across all sheets, as refreshed by the scheduled runner.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 7784.067
$ws.Range("I2").Value = 1209.0834
$ws.Range("K2").Value = 1209.0834
$ws.Range("M2").Value = -1096.0834
$ws.Range("H75").Value = 109212.14
$ws.Range("I75").Value = 191498.33
$ws.Range("J75").Value = 47497.5
$ws.Range("K75").Value = 191498.33
$ws.Range("L75").Value = 47497.5
$ws.Range("M75").Value = -190562.33
$ws.Range("N75").Value = -49369.5
$ws.Range("H78").Value = 109212.14
$ws.Range("I78").Value = 191498.33
$ws.Range("J78").Value = 47497.5
$ws.Range("K78").Value = 574494.99
$ws.Range("L78").Value = 142492.5
$ws.Range("M78").Value = -569814.99
$ws.Range("N78").Value = -151852.5
$ws.Range("H81").Value = 98747.75
$ws.Range("J81").Value = 98747.75
$ws.Range("L81").Value = 98747.75
$ws.Range("N81").Value = -100743.75
$ws.Range("H84").Value = 98747.75
$ws.Range("J84").Value = 98747.75
$ws.Range("L84").Value = 296243.25
$ws.Range("N84").Value = -306227.25
$ws.Range("H93").Value = 49000
$ws.Range("J93").Value = 49000
$ws.Range("L93").Value = 49000
$ws.Range("N93").Value = -53992
$ws.Range("H137").Value = 13895561
$ws.Range("I137").Value = 19232052
$ws.Range("K137").Value = 57696156
$ws.Range("M137").Value = -57693606

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 479.16666
$ws.Range("I5").Value = 525
$ws.Range("J5").Value = 456.25
$ws.Range("K5").Value = 525
$ws.Range("L5").Value = 456.25
$ws.Range("M5").Value = -413
$ws.Range("N5").Value = -680.25
$ws.Range("H39").Value = 32508
$ws.Range("I39").Value = 32508
$ws.Range("K39").Value = 32508
$ws.Range("M39").Value = -31988
$ws.Range("H50").Value = 386.5
$ws.Range("J50").Value = 367.83334
$ws.Range("L50").Value = 367.83334
$ws.Range("N50").Value = -1795.83334
$ws.Range("H74").Value = 307445.75
$ws.Range("I74").Value = 910646.6
$ws.Range("K74").Value = 910646.6
$ws.Range("M74").Value = -909772.6
$ws.Range("H77").Value = 307445.75
$ws.Range("I77").Value = 910646.6
$ws.Range("K77").Value = 4553233
$ws.Range("M77").Value = -4548865

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 479.16666
$ws.Range("I4").Value = 525
$ws.Range("J4").Value = 456.25
$ws.Range("K4").Value = 525
$ws.Range("L4").Value = 456.25
$ws.Range("M4").Value = -410
$ws.Range("N4").Value = -686.25
$ws.Range("H99").Value = 4137.7
$ws.Range("I99").Value = 3546
$ws.Range("J99").Value = 6504.5
$ws.Range("K99").Value = 3546
$ws.Range("L99").Value = 6504.5
$ws.Range("M99").Value = -2048
$ws.Range("N99").Value = -9500.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 48031.89
$ws.Range("I60").Value = 46498.4
$ws.Range("K60").Value = 46498.4
$ws.Range("M60").Value = -45987.4
$ws.Range("H62").Value = 10205.357
$ws.Range("I62").Value = 5297.5713
$ws.Range("K62").Value = 5297.5713
$ws.Range("M62").Value = -4673.5713
$ws.Range("H65").Value = 10205.357
$ws.Range("I65").Value = 5297.5713
$ws.Range("K65").Value = 26487.8565
$ws.Range("M65").Value = -23367.8565

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1360.2858
$ws.Range("J92").Value = 849.38464
$ws.Range("L92").Value = 2548.15392
$ws.Range("N92").Value = -5044.15392

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 494.33334
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 494.33334
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 494.33334
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -834.33334
$ws.Range("H31").Value = 1033.1666
$ws.Range("I31").Value = 819.8
$ws.Range("K31").Value = 819.8
$ws.Range("M31").Value = -527.8
$ws.Range("H37").Value = 1033.1666
$ws.Range("I37").Value = 819.8
$ws.Range("K37").Value = 819.8
$ws.Range("M37").Value = -542.8
$ws.Range("H70").Value = 13783.647
$ws.Range("I70").Value = 13611.625
$ws.Range("K70").Value = 13611.625
$ws.Range("M70").Value = -13341.625
$ws.Range("H73").Value = 13783.647
$ws.Range("I73").Value = 13611.625
$ws.Range("K73").Value = 13611.625
$ws.Range("M73").Value = -12675.625
$ws.Range("H80").Value = 2334.6667
$ws.Range("I80").Value = 2452.5
$ws.Range("K80").Value = 2452.5
$ws.Range("M80").Value = -1454.5
$ws.Range("H83").Value = 2334.6667
$ws.Range("I83").Value = 2452.5
$ws.Range("K83").Value = 12262.5
$ws.Range("M83").Value = -7270.5
$ws.Range("H96").Value = 34885.168
$ws.Range("J96").Value = 34885.168
$ws.Range("L96").Value = 34885.168
$ws.Range("N96").Value = -40377.168
$ws.Range("H102").Value = 1280.4546
$ws.Range("I102").Value = 969.8889
$ws.Range("K102").Value = 969.8889
$ws.Range("M102").Value = 652.1111
$ws.Range("H126").Value = 3389.875
$ws.Range("I126").Value = 2748.4
$ws.Range("K126").Value = 8245.200000000001
$ws.Range("M126").Value = -5775.200000000001
$ws.Range("H132").Value = 7571.7393
$ws.Range("I132").Value = 4837.8335
$ws.Range("K132").Value = 14513.5005
$ws.Range("M132").Value = -11983.5005

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 7251.3335
$ws.Range("J61").Value = 6800
$ws.Range("L61").Value = 6800
$ws.Range("N61").Value = -7204
$ws.Range("H113").Value = 7251.3335
$ws.Range("J113").Value = 6800
$ws.Range("L113").Value = 6800
$ws.Range("N113").Value = -11140

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 10159.308
$ws.Range("J52").Value = 39777
$ws.Range("L52").Value = 39777
$ws.Range("N52").Value = -40229
$ws.Range("H70").Value = 30361.666
$ws.Range("J70").Value = 27990
$ws.Range("L70").Value = 27990
$ws.Range("N70").Value = -28620
$ws.Range("H73").Value = 30361.666
$ws.Range("J73").Value = 27990
$ws.Range("L73").Value = 27990
$ws.Range("N73").Value = -30174
$ws.Range("H81").Value = 6095.769
$ws.Range("I81").Value = 1320.7142
$ws.Range("J81").Value = 11666.667
$ws.Range("K81").Value = 2641.4284
$ws.Range("L81").Value = 23333.334
$ws.Range("M81").Value = -1580.4284
$ws.Range("N81").Value = -25455.334
$ws.Range("H84").Value = 6095.769
$ws.Range("I84").Value = 1320.7142
$ws.Range("J84").Value = 11666.667
$ws.Range("K84").Value = 13207.142
$ws.Range("L84").Value = 116666.67
$ws.Range("M84").Value = -7903.142
$ws.Range("N84").Value = -127274.67
$ws.Range("H100").Value = 1474.0454
$ws.Range("I100").Value = 642.8823
$ws.Range("J100").Value = 4300
$ws.Range("K100").Value = 1285.7646
$ws.Range("L100").Value = 8600
$ws.Range("M100").Value = -744.7646
$ws.Range("N100").Value = -9682

Write-Host "Applied Zalera_Profits market data refresh"
